# Applies the "Daten Doppelfolge 6" edit to the "Nights" sheet:
#  - Adds a new row 7 of data (reusing existing names)
#  - Recolors a number of existing cells (fill colors change / get added)
#  - Updates the sheet's zoom level and active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nights")

# --- Color constants (Excel BGR-packed "Color" values) --------------------
$RED        = 255        # FFFF0000
$GREEN      = 5296274    # FF92D050
$YELLOW     = 65535       # FFFFFF00
$DARKGREEN  = 5287936    # FF00B050
$xlNone     = -4142

# --- Add new row 7 (Doppelfolge 6) -----------------------------------------
$ws.Range("B7").Value = "Nadja"
$ws.Range("C7").Value = "Anna"
$ws.Range("D7").Value = "Deisy"
$ws.Range("E7").Value = "Selina"
$ws.Range("F7").Value = "Nasti"
$ws.Range("G7").Value = "Sophia"
$ws.Range("H7").Value = "Ina"
$ws.Range("I7").Value = "Chiara"
$ws.Range("J7").Value = "Tori"
$ws.Range("K7").Value = "Camelia"
$ws.Range("L7").Value = 5

# --- Cells that stay/become red or green (pre-existing palette colors) ----
$ws.Range("B2").Interior.Color = $RED
$ws.Range("D2").Interior.Color = $GREEN
$ws.Range("E2").Interior.Color = $RED
$ws.Range("F2").Interior.Color = $RED
$ws.Range("K2").Interior.Color = $RED

$ws.Range("B3").Interior.Color = $RED
$ws.Range("D3").Interior.Color = $GREEN
$ws.Range("F3").Interior.Color = $RED
$ws.Range("K3").Interior.Color = $RED

$ws.Range("B4").Interior.Color = $RED
$ws.Range("D4").Interior.Color = $GREEN
$ws.Range("F4").Interior.Color = $RED

$ws.Range("B5").Interior.Color = $GREEN
$ws.Range("D5").Interior.Color = $GREEN
$ws.Range("F5").Interior.Color = $RED

$ws.Range("B6").Interior.Color = $GREEN
$ws.Range("D6").Interior.Color = $GREEN
$ws.Range("F6").Interior.Color = $RED
$ws.Range("G6").Interior.Color = $RED

$ws.Range("B7").Interior.Color = $GREEN
$ws.Range("C7").Interior.Color = $RED
$ws.Range("D7").Interior.Color = $GREEN
$ws.Range("G7").Interior.Color = $RED

# --- Cells whose fill gets removed entirely (formerly red) -----------------
$ws.Range("C5").Interior.Pattern = $xlNone
$ws.Range("K5").Interior.Pattern = $xlNone

# --- Cells that become yellow -----------------------------------------------
$ws.Range("C4").Interior.Color = $YELLOW
$ws.Range("E4").Interior.Color = $YELLOW
$ws.Range("H4").Interior.Color = $YELLOW
$ws.Range("I4").Interior.Color = $YELLOW
$ws.Range("J4").Interior.Color = $YELLOW
$ws.Range("K4").Interior.Color = $YELLOW
$ws.Range("H5").Interior.Color = $YELLOW
$ws.Range("I5").Interior.Color = $YELLOW
$ws.Range("J5").Interior.Color = $YELLOW
$ws.Range("C6").Interior.Color = $YELLOW
$ws.Range("E6").Interior.Color = $YELLOW
$ws.Range("I6").Interior.Color = $YELLOW
$ws.Range("J6").Interior.Color = $YELLOW
$ws.Range("K6").Interior.Color = $YELLOW
$ws.Range("H6").Interior.Color = $YELLOW

# --- Cells that become dark green -------------------------------------------
$ws.Range("G3").Interior.Color = $DARKGREEN
$ws.Range("G4").Interior.Color = $DARKGREEN
$ws.Range("G5").Interior.Color = $DARKGREEN

# --- View state: zoom + active cell selection -------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 142
$ws.Range("G7").Select()
